# The dataset (Repollo @ Terminal Hortofrutícola Agro Chillán) is a rolling
# weekly log. This edit inserts a new week's record at the top of the data
# block (row 59), which pushes every existing record down by one row, and
# the oldest record that falls off the bottom (old row 136) is re-appended
# as the new last row (137).
#
# Columns A,B,C,E,F,G,N,Q,R are constant for every data row in this block,
# so only D (Fecha), H (Variedad), I (Calidad), J (Volumen), K (Precio
# mínimo), L (Precio máximo), M (Precio promedio ponderado), O (Origen) and
# P (Precio $/Kg) need to move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D","H","I","J","K","L","M","O","P")

$firstDataRow = 59
$lastDataRow = 136
$newLastRow = 137

# 1) Capture the row that is about to fall off the bottom (old row 136) so
#    it can be re-inserted as the new row 137 once the shift is done.
$carry = @{}
foreach ($c in $cols) {
    $carry[$c] = $ws.Range("$c$lastDataRow").Value2
}
$carryA = $ws.Range("A$lastDataRow").Value2
$carryB = $ws.Range("B$lastDataRow").Value2
$carryC = $ws.Range("C$lastDataRow").Value2
$carryE = $ws.Range("E$lastDataRow").Value2
$carryF = $ws.Range("F$lastDataRow").Value2
$carryG = $ws.Range("G$lastDataRow").Value2
$carryN = $ws.Range("N$lastDataRow").Value2
$carryQ = $ws.Range("Q$lastDataRow").Value2
$carryR = $ws.Range("R$lastDataRow").Value2

# 2) Shift every existing record down by one row (bottom-up so we never
#    overwrite a source row before it has been read).
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $src = $r - 1
    foreach ($c in $cols) {
        $val = $ws.Range("$c$src").Value2
        $ws.Range("$c$r").Value = $val
    }
}

# 3) Write the new week's record into the now-vacated first data row.
$ws.Range("D$firstDataRow").Value = 44494
$ws.Range("H$firstDataRow").Value = "Crespo record"
$ws.Range("I$firstDataRow").Value = "Primera"
$ws.Range("J$firstDataRow").Value = 200
$ws.Range("K$firstDataRow").Value = 600
$ws.Range("L$firstDataRow").Value = 700
$ws.Range("M$firstDataRow").Value = 650
$ws.Range("O$firstDataRow").Value = "Región del Maule"
$ws.Range("P$firstDataRow").Value = 650

# 4) Re-append the record that fell off the bottom as the new last row.
$ws.Range("A$newLastRow").Value = $carryA
$ws.Range("B$newLastRow").Value = $carryB
$ws.Range("C$newLastRow").Value = $carryC
$ws.Range("D$newLastRow").NumberFormat = $ws.Range("D$lastDataRow").NumberFormat
$ws.Range("D$newLastRow").Value = $carry["D"]
$ws.Range("E$newLastRow").Value = $carryE
$ws.Range("F$newLastRow").Value = $carryF
$ws.Range("G$newLastRow").Value = $carryG
$ws.Range("H$newLastRow").Value = $carry["H"]
$ws.Range("I$newLastRow").Value = $carry["I"]
$ws.Range("J$newLastRow").Value = $carry["J"]
$ws.Range("K$newLastRow").Value = $carry["K"]
$ws.Range("L$newLastRow").Value = $carry["L"]
$ws.Range("M$newLastRow").Value = $carry["M"]
$ws.Range("N$newLastRow").Value = $carryN
$ws.Range("O$newLastRow").Value = $carry["O"]
$ws.Range("P$newLastRow").Value = $carry["P"]
$ws.Range("Q$newLastRow").Value = $carryQ
$ws.Range("R$newLastRow").Value = $carryR
